$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new article ("SVG 05 - Path 2") is inserted as the newest entry at the top
# of the table, pushing every existing article row down by one.
# Insert a new row at position 2; this is a pure structural shift, so rows that
# used to be 2..24 become 3..25 with all of their cell values/styles untouched.
$ws.Rows.Item(2).Insert()

# Seed the new row 2 from row 3 (the row that used to be row 2) so it inherits
# the correct cell style (border formatting) used throughout the table, and the
# correct "tag-web" tag in column A (same tag as the row below it).
$ws.Range("A3:E3").Copy($ws.Range("A2:E2"))

# Fill in the new row with the new article's data. Values are entered in the
# same order as the original authoring (title, img, date, then site url) so the
# freshly created shared-string entries line up with the target workbook.
$ws.Cells.Item(2, 2).Value = "SVG 研究之路 (5) - Path 進階篇"
$ws.Cells.Item(2, 4).Value = "/img/articles/201406/20140612_1_01.jpg"
$ws.Cells.Item(2, 5).Value = "JUN 12TH, 2014"
$ws.Cells.Item(2, 3).Value = "/articles/201406/svg-05-path-2.html"

# Match the author's final selection (GA check around the new top row).
$ws.Range("C2").Select()
